$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values, repulled from source data
$updates = @{
    2  = 0
    6  = 4
    12 = -1
    13 = 3
    16 = 8
    20 = 0
    21 = 0
    26 = -3
    30 = 1
    31 = 5
    34 = 2
    35 = 0
    36 = -3
    38 = 0
    39 = -4
    42 = 3
    43 = 7
    44 = 3
    45 = -3
    46 = 1
    48 = 4
    51 = 3
    52 = -8
    53 = -1
    54 = 1
    55 = -1
    56 = 8
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
